# Generate Report for Handback
# The f889c0fd-9ec9-4624-b0c8-bc6f7f90dd34 file has now been handed back and is
# in sync with en-US, so the localization-status report is refreshed to reflect
# that: status flips from "Ready for handoff" -> "Handed back: in sync with
# en-US" everywhere it is shown, the per-language "Latest Handback DateTime" is
# stamped with the new handback time, and the stale "handback file is not the
# latest" error is cleared now that the handback succeeded.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: row for f889c0fd...md, zh-cn (E3) and de-de (F3) columns ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

# --- zh-cn sheet: row for f889c0fd...md (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusText
$zhcn.Range("K3").Value = "2016-08-23 18:51:59"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).ColumnWidth = 13.7470528738839

# --- de-de sheet: row for f889c0fd...md (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusText
$dede.Range("K3").Value = "2016-08-23 18:52:18"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).ColumnWidth = 13.7470528738839
